{"js": "// Insert a new paragraph \"dotnet ef migrations add InitialData\" right\n// after the existing \"dotnet ef migrations add ColumnPesoCategoria\"\n// paragraph, matching that paragraph's formatting (bold, 14pt / sz 28).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"dotnet ef migrations add ColumnPesoCategoria\") !== -1) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(\"Could not locate the 'ColumnPesoCategoria' migration paragraph.\");\n}\n\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"dotnet ef migrations add InitialData\",\n  \"After\"\n);\n\n// Make sure the new run/paragraph carries the same bold + 14pt (sz 28)\n// formatting as the paragraph it follows (inherited already, set\n// explicitly too for robustness).\nnewParagraph.font.bold = true;\nnewParagraph.font.size = 14;\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"dotnet ef migrations add InitialData\" right\n# after the existing \"dotnet ef migrations add ColumnPesoCategoria\"\n# paragraph, matching that paragraph's formatting (bold, 14pt / sz 28).\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = \"dotnet ef migrations add ColumnPesoCategoria\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate the 'ColumnPesoCategoria' migration paragraph.\"\n}\n\n# $searchRange now spans the found text; Paragraphs(1) is its paragraph.\n$anchorParagraph = $searchRange.Paragraphs(1)\n\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph = $anchorParagraph.Next()\n$newParagraph.Range.Text = \"dotnet ef migrations add InitialData\"\n\n# Make sure the new run/paragraph carries the same bold + 14pt (sz 28)\n# formatting as the paragraph it follows (inherited already, set\n# explicitly too for robustness).\n$newParagraph.Range.Font.Bold = $true\n$newParagraph.Range.Font.Size = 14\n"}
